$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# Expand the "Tableau4" table by one row; this grows the table ref/autofilter
# and brings along the calculated column definition.
$newRow = $lo.ListRows.Add()

# Copy formatting (styles/borders) from the last existing data row (13) down
# to the freshly added row (14).
$ws.Range("A13:H13").Copy($ws.Range("A14:H14"))

# The Copy above only pasted the literal value of D13's formula; restore the
# structured-reference formula for the calculated column in the new row.
$ws.Range("D14").Formula = "=Tableau4[[#This Row],[Heure fin]]-Tableau4[[#This Row],[Heure début]]"

# Complete the previously unfinished entry in row 13 (end time, place, theme,
# description).
$ws.Range("C13").Value = 0.39861111111111108
$ws.Range("E13").Value = "CPNV"
$ws.Range("F13").Value = "Cahier des Charges"
$ws.Range("G13").Value = "J'ai fait une mise en comum du cahier des charges avec Evann"

# Fill in the new journal entry on row 14.
$ws.Range("A14").Value = 45056
$ws.Range("B14").Value = 0.64097222222222217
$ws.Range("C14").Value = 0.70277777777777783
$ws.Range("E14").Value = "CPNV"
$ws.Range("F14").Value = "Base de données"
$ws.Range("G14").Value = "J'ai commencé le script pour créer la base de données"

# Row 3 reverts to the default (15pt) but explicit custom height.
$ws.Rows.Item(3).RowHeight = 15

# Leave the selection on the last cell of the new row, like the author did.
$null = $ws.Range("H14").Select()
